$p = $ppt.ActivePresentation

# Duplicate slide 5 (the "Type 5" prototype slide) to create the new slide 6,
# then move it to the end of the deck and retitle it "Type 6".
$src = $p.Slides.Item(5)
$dup = $src.Duplicate()
$dup.MoveTo($p.Slides.Count)

$newSlide = $p.Slides.Item($p.Slides.Count)
$newSlide.Shapes.Item(1).TextFrame.TextRange.Text = "Type 6"
